$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 228, shifting the rest of the
# data block (old rows 228-266) down to rows 230-268.
$ws.Rows("228:229").Insert()

# Row 228: new weekly entry for "Magnum" variety.
$ws.Cells.Item(228, 1).Value = 2
$ws.Cells.Item(228, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(228, 3).Value = "Coquimbo"
$ws.Cells.Item(228, 4).Value = 45015
$ws.Cells.Item(228, 5).Value = 4
$ws.Cells.Item(228, 6).Value = 100112031
$ws.Cells.Item(228, 7).Value = "Poroto verde"
$ws.Cells.Item(228, 8).Value = "Magnum"
$ws.Cells.Item(228, 9).Value = "Primera"
$ws.Cells.Item(228, 10).Value = 520
$ws.Cells.Item(228, 11).Value = 14000
$ws.Cells.Item(228, 12).Value = 15000
$ws.Cells.Item(228, 13).Value = 14500
$ws.Cells.Item(228, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(228, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(228, 16).Value = 580
$ws.Cells.Item(228, 17).Value = 25
$ws.Cells.Item(228, 18).Value = "Hortaliza"

# Row 229: new weekly entry for "Sin especificar" variety.
$ws.Cells.Item(229, 1).Value = 2
$ws.Cells.Item(229, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(229, 3).Value = "Coquimbo"
$ws.Cells.Item(229, 4).Value = 45015
$ws.Cells.Item(229, 5).Value = 4
$ws.Cells.Item(229, 6).Value = 100112031
$ws.Cells.Item(229, 7).Value = "Poroto verde"
$ws.Cells.Item(229, 8).Value = "Sin especificar"
$ws.Cells.Item(229, 9).Value = "Primera"
$ws.Cells.Item(229, 10).Value = 600
$ws.Cells.Item(229, 11).Value = 22000
$ws.Cells.Item(229, 12).Value = 23000
$ws.Cells.Item(229, 13).Value = 22500
$ws.Cells.Item(229, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(229, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(229, 16).Value = 900
$ws.Cells.Item(229, 17).Value = 25
$ws.Cells.Item(229, 18).Value = "Hortaliza"
